$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '49.588.59'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').Value = '2.636.34'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '112.94'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.59%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '323.97'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.09%  '
$ws.Range('E7').Value = '  -0.55%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.545'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.25%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.80'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.96%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '19.83'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.63%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0812'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('E13').Value = '  +1.14%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.31'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.52%  '
$ws.Range('D15').Value = '3.050.23'
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('D16').Value = '2.640.99'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('E17').Value = '  -2.11%  '
$ws.Range('D18').Value = '49.556.96'
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('E19').Value = '  -3.05%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.89'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.19%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.70'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.66%  '
$ws.Range('D22').Value = '0.0₃0947'
$ws.Range('E22').Value = '  -1.28%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '270.12'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -4.09%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '69.01'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -5.45%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.55'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.21%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '26.33'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.11%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.32'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +4.11%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.22'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.77%  '
$ws.Range('E30').Value = '  -3.00%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '35.01'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -4.09%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '49.62'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.18%  '
$ws.Range('E33').Value = '  +1.15%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0814'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.52%  '
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '18.97'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -3.55%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.93'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +4.20%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.04'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.50%  '
$ws.Range('E39').Value = '  +0.45%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '126.93'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.92%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '22.54'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.76%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.111'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.31%  '
$ws.Range('E43').Value = '  +3.62%  '
$ws.Range('E44').Value = '  -3.11%  '
$ws.Range('D45').Value = '2.060.21'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('E46').Value = '  -3.32%  '
$ws.Range('E47').Value = '  +7.13%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.16'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -6.35%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.94'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.18%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '5.23'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.89%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '59.06'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.59%  '
